$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G holds "K" values (K = strike); regen save_data to use K instead of Strike#
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 3
$ws.Range("G7").Value = 0
